# Update the "Data de emissão" (issue date / time) values on the active sheet.
# B2 holds the date string, C2 holds the time string. Both are plain text
# (shared strings), so we assign them explicitly as text to avoid Excel
# auto-converting them to date/time serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "30/08/2022"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "08:44:47"
